# Updated cryptos list values (Price / Volume(1h)) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.949.13'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '2.355.14'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('E5').Value = '  +5.33%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '241.92'
$ws.Range('E6').Value = '  +3.10%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '76.60'
$ws.Range('E7').Value = '  +3.40%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.640'
$ws.Range('E9').Value = '  +21.06%  '
$ws.Range('E10').Value = '  +3.83%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '57.41'
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '33.47'
$ws.Range('E12').Value = '  +22.32%  '
$ws.Range('E13').Value = '  +13.05%  '
$ws.Range('E14').Value = '  +1.48%  '
$ws.Range('D15').Value = '2.705.44'
$ws.Range('E15').Value = '  -0.16%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '16.83'
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.932'
$ws.Range('E17').Value = '  +5.90%  '
$ws.Range('D18').Value = '2.357.14'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').Value = '43.863.24'
$ws.Range('E19').Value = '  +0.39%  '
$ws.Range('E20').Value = '  +2.06%  '
$ws.Range('E21').Value = '  +3.01%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '77.73'
$ws.Range('E22').Value = '  +2.74%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '262.69'
$ws.Range('E23').Value = '  +4.55%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('E25').Value = '  +2.15%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.63'
$ws.Range('E26').Value = '  -5.38%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.98'
$ws.Range('E27').Value = '  +7.39%  '
$ws.Range('E28').Value = '  +16.80%  '
$ws.Range('E29').Value = '  +2.41%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '23.14'
$ws.Range('E30').Value = '  +2.70%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '175.27'
$ws.Range('E31').Value = '  +1.84%  '
$ws.Range('E32').Value = '  -3.65%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.137'
$ws.Range('E33').Value = '  +4.47%  '
$ws.Range('E34').Value = '  +5.78%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0764'
$ws.Range('E35').Value = '  +8.68%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.42'
$ws.Range('E36').Value = '  +5.64%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.81'
$ws.Range('E37').Value = '  +1.47%  '
$ws.Range('E38').Value = '  -0.51%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.42'
$ws.Range('E39').Value = '  -2.86%  '
$ws.Range('E40').Value = '  +7.09%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.217'
$ws.Range('E41').Value = '  +21.93%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '19.30'
$ws.Range('E42').Value = '  -1.14%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '9.17'
$ws.Range('E43').Value = '  +2.91%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.107'
$ws.Range('E44').Value = '  +10.10%  '
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.53'
$ws.Range('E46').Value = '  +10.75%  '
$ws.Range('E47').Value = '  +4.10%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '102.09'
$ws.Range('E49').Value = '  +1.49%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.58'
$ws.Range('E50').Value = '  +3.14%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '56.24'
$ws.Range('E51').Value = '  +10.26%  '
